$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a brand-new row at 279; this shifts the existing rows 279-297
# down to 280-298 (and the sheet dimension grows from R297 to R298).
$ws.Rows.Item(279).Insert()

# Populate the newly inserted row 279 with its data. Columns A,B,C,E,F,
# G,H,I,R are identical to the surrounding rows for this market/product,
# so copy them straight from row 280 (the row that used to be 279).
$ws.Range("A279").Value = $ws.Range("A280").Value()
$ws.Range("B279").Value = $ws.Range("B280").Value()
$ws.Range("C279").Value = $ws.Range("C280").Value()
$ws.Range("D279").Value = 45021
$ws.Range("E279").Value = $ws.Range("E280").Value()
$ws.Range("F279").Value = $ws.Range("F280").Value()
$ws.Range("G279").Value = $ws.Range("G280").Value()
$ws.Range("H279").Value = $ws.Range("H280").Value()
$ws.Range("I279").Value = $ws.Range("I280").Value()
$ws.Range("J279").Value = 140
$ws.Range("K279").Value = 6000
$ws.Range("L279").Value = 7000
$ws.Range("M279").Value = 6571
$ws.Range("N279").Value = "$/caja 50 unidades"
$ws.Range("O279").Value = "Región del Maule"
$ws.Range("P279").Value = 131
$ws.Range("Q279").Value = 50
$ws.Range("R279").Value = $ws.Range("R280").Value()
